$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column A entirely; remaining columns (B:F) shift left to become A:E
$ws.Range("A:A").Delete()
